$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the URL text into A1 with wrap text enabled
$ws.Range("A1").Value = "https://docs.google.com/spreadsheets/d/1NdRKcFUBwytjB3xhmOh6yPCHCwd7VqLOaN6Vfuk-G4w/edit?gid=276038343#gid=276038343"
$ws.Range("A1").WrapText = $true

# Widen column A (closest achievable value to 40.5703125 given engine's
# internal half-point column width storage / pixel rounding)
$ws.Columns.Item(1).ColumnWidth = 39.66

# Clear the selection stored in the sheet view (select A1 instead of G10)
$ws.Range("A1").Select()
